# The commit swaps the contents of ppt/theme/theme1.xml (stock "Office
# Theme") and ppt/theme/theme2.xml (custom "Integral" theme) so that the
# theme actually driving the slide master/slides becomes the plain
# "Office Theme" color palette (the "Integral" palette moves to the
# otherwise-unused theme1.xml slot, which - in this package - is only
# wired to the Notes Master).
#
# The live Presentation/SlideMaster/Slide object model here always
# resolves to the single theme part that backs the slide master (i.e.
# the package's theme2.xml); that is the only theme part reachable
# through PowerPoint's object model, so we reproduce the net effect of
# the swap for every property the OM exposes: replace the 12 theme
# colors with the stock "Office" palette, via ThemeColorScheme so the
# <a:clrScheme>/<a:theme> name attributes are left untouched by the
# write (matches how PowerPoint itself only rewrites the color nodes
# when you tweak theme colors one-by-one).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# msoThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink.
# Target values are the stock Office theme RGB colors; RGB must be
# supplied VBA-style (R | G<<8 | B<<16).
function HexToVbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = HexToVbaRgb $officeThemeHex[$i - 1]
}
